# Appends " (UNL)" as a new, separate run immediately after the given
# search phrase (which is assumed to end at the end of its paragraph).
# A genuine sibling <w:r> is produced (rather than the text being merged
# into the preceding run) by toggling a run-level property on the newly
# inserted text: Word's run-normalizer only fuses adjacent runs that are
# formatting-identical, so briefly flipping Bold on/off on just the new
# span keeps it a distinct run while leaving no visible formatting change.
function Add-UnlSuffix {
    param(
        [string]$SearchText
    )

    $d = $word.ActiveDocument
    $rng = $d.Content
    $found = $rng.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $SearchText"
    }

    $insertStart = $rng.End
    $rng.InsertAfter(" (UNL)")

    $newRng = $d.Range($insertStart, $insertStart + 6)
    $newRng.Bold = $true
    $newRng.Bold = $false
}

Add-UnlSuffix "University of Nebraska-Lincoln"
Add-UnlSuffix "Karl M. Kuntzelman, Ph.D."
Add-UnlSuffix "Michael D. Dodd, Ph.D."
Add-UnlSuffix "Matthew R. Johnson, Ph.D."
